$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 843   # was 841
$ws.Range("F4").Value = 805   # was 802
$ws.Range("F6").Value = 1040   # was 1020
$ws.Range("F7").Value = 1102   # was 1086
$ws.Range("F9").Value = 150   # was 149
$ws.Range("F10").Value = 499   # was 494
$ws.Range("F11").Value = 233   # was 212
$ws.Range("F12").Value = 45   # was 42
$ws.Range("F13").Value = 1196   # was 1187
$ws.Range("F14").Value = 26580   # was 26421
$ws.Range("F15").Value = 2947   # was 2888
$ws.Range("F16").Value = 19   # was 18
$ws.Range("F17").Value = 210   # was 206
$ws.Range("F18").Value = 409   # was 403
$ws.Range("F19").Value = 39   # was 38
$ws.Range("F20").Value = 261   # was 248
$ws.Range("F21").Value = 507   # was 495
$ws.Range("F22").Value = 245   # was 244
$ws.Range("F23").Value = 201   # was 194
$ws.Range("F24").Value = 312   # was 307
$ws.Range("F25").Value = 17   # was 16
$ws.Range("F26").Value = 621   # was 616
$ws.Range("F27").Value = 160   # was 151
$ws.Range("F28").Value = 64   # was 61
$ws.Range("F29").Value = 444   # was 443
$ws.Range("F30").Value = 43   # was 39
$ws.Range("F32").Value = 545   # was 535
$ws.Range("F33").Value = 217   # was 215

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 338   # was 332
$ws.Range("F7").Value = 617   # was 570
$ws.Range("F8").Value = 69   # was 64
$ws.Range("F9").Value = 256   # was 255
$ws.Range("F10").Value = 4194   # was 4181
$ws.Range("F17").Value = 31   # was 32
$ws.Range("F18").Value = 41   # was 40

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 224   # was 220
$ws.Range("F4").Value = 1057   # was 1033
$ws.Range("F5").Value = 277   # was 276

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 224   # was 220
$ws.Range("F4").Value = 1057   # was 1033
$ws.Range("F6").Value = 805   # was 802
$ws.Range("F11").Value = 338   # was 332
$ws.Range("F12").Value = 277   # was 276
$ws.Range("F13").Value = 617   # was 570
$ws.Range("F14").Value = 1041   # was 1020
$ws.Range("F15").Value = 1102   # was 1086
$ws.Range("F16").Value = 150   # was 149
$ws.Range("F17").Value = 499   # was 494
$ws.Range("F18").Value = 233   # was 212
$ws.Range("F19").Value = 45   # was 42
$ws.Range("F20").Value = 1196   # was 1187
$ws.Range("F21").Value = 69   # was 64
$ws.Range("F22").Value = 256   # was 255
$ws.Range("F27").Value = 2947   # was 2888
$ws.Range("F28").Value = 210   # was 206
$ws.Range("F31").Value = 409   # was 403
$ws.Range("F33").Value = 31   # was 32
$ws.Range("F34").Value = 261   # was 248
$ws.Range("F35").Value = 507   # was 495
$ws.Range("F36").Value = 245   # was 244
$ws.Range("F37").Value = 313   # was 307
$ws.Range("F38").Value = 17   # was 16
$ws.Range("F39").Value = 621   # was 616
$ws.Range("F40").Value = 41   # was 40
$ws.Range("F41").Value = 160   # was 151
$ws.Range("F42").Value = 64   # was 61
$ws.Range("F45").Value = 43   # was 39
$ws.Range("F47").Value = 545   # was 535
$ws.Range("F48").Value = 217   # was 215
